# Update "想去人数" (want-to-go count, column F) figures to the freshly
# scraped values, as published to gh-pages at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 4598
$ws.Range("F3").Value  = 2699
$ws.Range("F5").Value  = 2696
$ws.Range("F9").Value  = 1703
$ws.Range("F10").Value = 725
$ws.Range("F11").Value = 464
$ws.Range("F12").Value = 177
$ws.Range("F14").Value = 43
$ws.Range("F16").Value = 85
$ws.Range("F21").Value = 632
$ws.Range("F22").Value = 730
$ws.Range("F25").Value = 485
$ws.Range("F27").Value = 1385
$ws.Range("F28").Value = 285
$ws.Range("F30").Value = 1368
$ws.Range("F31").Value = 2229
$ws.Range("F37").Value = 86
$ws.Range("F38").Value = 748
$ws.Range("F39").Value = 1427
$ws.Range("F40").Value = 178
$ws.Range("F42").Value = 471
$ws.Range("F44").Value = 102

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 17

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 4598
$ws.Range("F3").Value  = 2699
$ws.Range("F4").Value  = 2696
$ws.Range("F5").Value  = 1703
$ws.Range("F8").Value  = 725
$ws.Range("F9").Value  = 464
$ws.Range("F10").Value = 177
$ws.Range("F12").Value = 46
$ws.Range("F14").Value = 85
$ws.Range("F18").Value = 632
$ws.Range("F19").Value = 730
$ws.Range("F25").Value = 485
$ws.Range("F27").Value = 1385
$ws.Range("F28").Value = 285
$ws.Range("F32").Value = 2229
$ws.Range("F37").Value = 17
$ws.Range("F41").Value = 86
$ws.Range("F42").Value = 748
$ws.Range("F43").Value = 1427
$ws.Range("F45").Value = 178
$ws.Range("F46").Value = 471
$ws.Range("F48").Value = 102
